$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "64.396.12"
$ws.Range("E2").Value = "  +0.21%  "
$ws.Range("D3").Value = "3.518.90"
$ws.Range("E3").Value = "  +0.51%  "
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").Value = "'592.16"
$ws.Range("E5").Value = "  +1.44%  "
$ws.Range("D6").Value = "'134.84"
$ws.Range("E6").Value = "  -0.05%  "
$ws.Range("E8").Value = "  +0.33%  "
$ws.Range("D9").Value = "'7.58"
$ws.Range("E9").Value = "  +6.11%  "
$ws.Range("E10").Value = "  +0.37%  "
$ws.Range("E11").Value = "  +3.91%  "
$ws.Range("D12").Value = "4.117.34"
$ws.Range("E12").Value = "  +0.47%  "
$ws.Range("E13").Value = "  +1.59%  "
$ws.Range("E14").Value = "  +1.07%  "
$ws.Range("D15").Value = "3.516.86"
$ws.Range("E15").Value = "  +0.47%  "
$ws.Range("D16").Value = "'25.89"
$ws.Range("E16").Value = "  -2.07%  "
$ws.Range("D17").Value = "64.379.26"
$ws.Range("E17").Value = "  +0.16%  "
$ws.Range("D18").Value = "'9.98"
$ws.Range("E18").Value = "  +2.38%  "
$ws.Range("E19").Value = "  +3.47%  "
$ws.Range("D20").Value = "'13.60"
$ws.Range("E20").Value = "  -1.70%  "
$ws.Range("D21").Value = "'394.62"
$ws.Range("E21").Value = "  +2.78%  "
$ws.Range("D22").Value = "'0.578"
$ws.Range("E22").Value = "  +1.54%  "
$ws.Range("D23").Value = "3.659.30"
$ws.Range("E23").Value = "  +0.50%  "
$ws.Range("D24").Value = "'74.78"
$ws.Range("E24").Value = "  +1.22%  "
$ws.Range("E25").Value = "  -0.04%  "
$ws.Range("D26").Value = "'5.73"
$ws.Range("E26").Value = "  +0.19%  "
$ws.Range("E27").Value = "  +3.48%  "
$ws.Range("B28").Value = "Binance-PegBSC-USD"
$ws.Range("C28").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D28").Value = "'1.00"
$ws.Range("E28").Value = "  -0.06%  "
$ws.Range("B29").Value = "RenderToken"
$ws.Range("C29").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D29").Value = "'7.41"
$ws.Range("E29").Value = "  -1.73%  "
$ws.Range("D30").Value = "'2.26"
$ws.Range("E30").Value = "  +1.77%  "
$ws.Range("D31").Value = "'8.32"
$ws.Range("E31").Value = "  +0.25%  "
$ws.Range("D32").Value = "'1.47"
$ws.Range("E32").Value = "  -6.00%  "
$ws.Range("E33").Value = "  +7.92%  "
$ws.Range("D34").Value = "3.549.21"
$ws.Range("E34").Value = "  +0.82%  "
$ws.Range("E35").Value = "  +0.05%  "
$ws.Range("D36").Value = "'23.46"
$ws.Range("E36").Value = "  -0.73%  "
$ws.Range("D37").Value = "'5.38"
$ws.Range("E37").Value = "  +1.17%  "
$ws.Range("E38").Value = "  +1.74%  "
$ws.Range("E39").Value = "  +1.90%  "
$ws.Range("D40").Value = "'167.35"
$ws.Range("D41").Value = "'0.0793"
$ws.Range("E41").Value = "  +1.39%  "
$ws.Range("E42").Value = "  +0.39%  "
$ws.Range("D43").Value = "'25.65"
$ws.Range("E43").Value = "  -1.72%  "
$ws.Range("E44").Value = "  -0.03%  "
$ws.Range("D45").Value = "'4.46"
$ws.Range("E45").Value = "  +1.03%  "
$ws.Range("D46").Value = "'1.67"
$ws.Range("E46").Value = "  +2.96%  "
$ws.Range("E47").Value = "  -1.89%  "
$ws.Range("D48").Value = "'6.81"
$ws.Range("E48").Value = "  +0.73%  "
$ws.Range("D49").Value = "2.412.77"
$ws.Range("E49").Value = "  -2.25%  "
$ws.Range("D50").Value = "'0.899"
$ws.Range("E50").Value = "  -2.25%  "
$ws.Range("E51").Value = "  +0.09%  "
